$d = $word.ActiveDocument

# 1. Update the main title text
$d.Content.Find.Execute("AI Conversational Revolution", $true, $false, $false, $false, $false, $true, 1, $false, "Revolutionizing Artificial Intelligence with ChatGPT", 2)

# 2. Append the large block of new paragraphs (sections + closing Gmail-tool confirmation
#    paragraphs) after the title, using a raw WordOpenXML package so every paragraph gets
#    exactly the pPr/rPr formatting captured in the target document.
$rng = $d.Content
$rng.Collapse(0)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Introductory Title</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ChatGPT is a revolutionary artificial intelligence language model developed by OpenAI. It has taken the world by storm with its ability to understand and respond to human input in a remarkably natural way. One of the key features of ChatGPT is its ability to engage in conversation, using context and understanding to generate human-like responses.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Learning and Improvement</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ChatGPT''s ability to learn and improve over time is one of its most impressive aspects. Through machine learning algorithms and vast amounts of training data, ChatGPT is able to refine its understanding of language and generate more accurate and informative responses. This has made it an invaluable tool for a wide range of applications, from customer service to language translation.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Augmenting Human Capabilities</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Another important aspect of ChatGPT is its potential to augment human capabilities. By automating tasks such as data entry and customer support, ChatGPT can free up humans to focus on more complex and creative tasks. Additionally, ChatGPT''s ability to process and analyze large amounts of data makes it an invaluable tool for researchers and scientists.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Risks and Challenges</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Despite its many benefits, ChatGPT also raises important questions about the future of work and the potential risks of relying on artificial intelligence. As ChatGPT and other AI models become increasingly sophisticated, there is a risk that they will displace human workers and exacerbate existing social and economic inequalities.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Transformative Impact</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Ultimately, the true potential of ChatGPT will depend on how it is used and developed in the years to come. If harnessed responsibly and used to augment human capabilities, ChatGPT could have a transformative impact on a wide range of industries and fields.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Far-Reaching Implications</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Additionally, ChatGPT''s capabilities have far-reaching implications for numerous fields and careers. For instance, in healthcare, ChatGPT can assist with patient data analysis, medical research, and even provide emotional support to patients. In education, ChatGPT can help with personalized learning, grading, and even create customized educational content. Additionally, ChatGPT can aid in language translation, facilitating global communication and collaboration. Its potential applications are vast, and its impact will be felt across multiple industries.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve">The Power of ChatGPT</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Revolutionizing Language Models</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ChatGPT is a revolutionary artificial intelligence language model developed by OpenAI. It has taken the world by storm with its ability to understand and respond to human input in a remarkably natural way.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">One of the key features of ChatGPT is its ability to engage in conversation, using context and understanding to generate human-like responses.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">ChatGPT''s Ability to Learn and Improve</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ChatGPT''s ability to learn and improve over time is one of its most impressive aspects. Through machine learning algorithms and vast amounts of training data, ChatGPT is able to refine its understanding of language and generate more accurate and informative responses.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Augmenting Human Capabilities</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Another important aspect of ChatGPT is its potential to augment human capabilities. By automating tasks such as data entry and customer support, ChatGPT can free up humans to focus on more complex and creative tasks.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Risks and Concerns</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Despite its many benefits, ChatGPT also raises important questions about the future of work and the potential risks of relying on artificial intelligence.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Transformative Impact</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Ultimately, the true potential of ChatGPT will depend on how it is used and developed in the years to come. If harnessed responsibly and used to augment human capabilities, ChatGPT could have a transformative impact on a wide range of industries and fields.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Far-Reaching Implications</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Additionally, ChatGPT''s capabilities have far-reaching implications for numerous fields and careers.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">In healthcare, ChatGPT can assist with patient data analysis, medical research, and even provide emotional support to patients. In education, ChatGPT can help with personalized learning, grading, and even create customized educational content.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Moreover, in computer careers, ChatGPT can assist with coding, debugging, and even generating new code.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve">Gmail and Write File Tool Confirmation</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">I''d be happy to send this document in an email. Before I do, I want to confirm that I''ll be using the Gmail tool.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Please make sure to enter your credentials in the Doc Manager extension.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">Additionally, I''ve used the write file tool previously. Are you sure you want to use it again? If so, I''ll proceed with caution.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve">I will use the Gmail tool.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve">Confirming Gmail Tool Use</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">I''ll be happy to send this document in an email. Before I do, I want to confirm that I''ll be using the Gmail tool.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Please make sure to enter your credentials in the Doc Manager extension.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Additionally, I''ve used the write file tool previously. Are you sure you want to use it again?</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b w:val="0"/><w:i w:val="0"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">If so, I''ll proceed with caution.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i w:val="0"/><w:color w:val="0000FF"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">I will use the Gmail tool.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
